# Add 2022-Q3 data: new sheet + new summary row on "总计"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new row on the "总计" (totals) sheet for the 2022-Q3 period,
#    right above the existing 2022-Q2 row, and renumber the index column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# Row-insert copies the formatting of the row above into B2:D2 - the
# original sheet keeps those cells unstyled, so strip it back off.
$total.Range("B2:D2").ClearFormats()

# A2 should carry the same index-column style as the rest of column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.03

# Renumber the rest of the index column (it used to run 0..3, now 1..4).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# 2) Insert a brand-new worksheet named "2022-Q3" right before "2022-Q2"
#    and populate it with the quarter's fund holdings.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Match the bold/bordered header style used on every other sheet.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A2 uses the same index-column style as the other sheets' leading column.
$total.Range("A3").Copy()
$q3.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$q3.Range("A2").Value = 0

# Columns B-G on this sheet are stored as text (fund codes need their
# leading zeros, and the numeric-looking figures are plain text too, just
# like on the other quarterly sheets) - force text format before writing,
# then drop the now-unneeded number-format style so the cells end up
# unstyled, same as on the sibling sheets.
$q3.Range("B2:G2").NumberFormat = "@"
$q3.Range("B2").Value = "001972"
$q3.Range("C2").Value = "前海开源沪港深智慧生活优选灵活配置混合"
$q3.Range("D2").Value = "0.56"
$q3.Range("E2").Value = "91.51"
$q3.Range("F2").Value = "5.38"
$q3.Range("G2").Value = "0.0301"
$q3.Range("B2:G2").ClearFormats()

$q3.Range("H2").Value = 9

# Restore the original active tab ("总计") now that the new sheet is no
# longer the freshly-created/active one.
$total.Activate()
